$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Collected Amount *" column (column G) is no longer populated/needed,
# so remove it entirely - this shifts Status/Verified/Folio No columns left
# (H->G, I->H, J->I) and drops the now-unused "Collected Amount *" shared
# string, matching a user selecting the whole column and deleting it.
$ws.Range("G1").EntireColumn.Select()
$ws.Range("G1").EntireColumn.Delete()
